$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column T (Pension Percepciones por persona) changes from a "measure" to a "dimension"
$ws.Range("T3").Value = "iaest-dimension:pension-percepciones-por-persona"
$ws.Range("T4").Value = "dim"
$ws.Range("T5").Value = "skos:Concept"

# Column W (CCAA) changes from iaest-measure:ccaa to sdmx-dimension:refArea
$ws.Range("W3").Value = "sdmx-dimension:refArea"
$ws.Range("W4").Value = "dim"
$ws.Range("W5").Value = "URI-Comunidad"

# New row 6 with reference to mapping file, using same formatting as row 5
$ws.Range("T5").Copy()
$ws.Range("T6").PasteSpecial(-4122)
$ws.Range("T6").Value = "mapping-pension-percepciones-por-persona.xlsx"
